$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.31250333333333
$ws.Range("H2").Value = 36.93751
$ws.Range("I2").Value = 0.6498350963072504
$ws.Range("J2").Value = 0.6498350963072506
$ws.Range("M2").Value = 4.889237666666667
$ws.Range("N2").Value = 14.667713
$ws.Range("O2").Value = 0.09529921759032918
$ws.Range("P2").Value = 0.09529921759032917
$ws.Range("Q2").Value = 60.19875506829223
$ws.Range("R2").Value = 541.7887956146301
$ws.Range("S2").Value = 0.06192877624081718
$ws.Range("T2").Value = 0.06192877624081718
$ws.Range("G3").Value = 12.31250333333333
$ws.Range("H3").Value = 36.93751
$ws.Range("I3").Value = 0.6498350963072504
$ws.Range("J3").Value = 0.6498350963072506
$ws.Range("O3").Value = 0.6664673019309815
$ws.Range("P3").Value = 0.6664673019309812
$ws.Range("Q3").Value = 420.9950814332822
$ws.Range("R3").Value = 3788.95573289954
$ws.Range("S3").Value = 0.4330938433359527
$ws.Range("T3").Value = 0.4330938433359526
$ws.Range("G4").Value = 12.31250333333333
$ws.Range("H4").Value = 36.93751
$ws.Range("I4").Value = 0.6498350963072504
$ws.Range("J4").Value = 0.6498350963072506
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06428533333333333
$ws.Range("N4").Value = 0.192856
$ws.Range("O4").Value = 0.001253026010776221
$ws.Range("P4").Value = 0.001253026010776221
$ws.Range("Q4").Value = 0.7915133809511111
$ws.Range("R4").Value = 7.123620428560001
$ws.Range("S4").Value = 0.0008142602783882557
$ws.Range("T4").Value = 0.0008142602783882557
$ws.Range("G5").Value = 12.31250333333333
$ws.Range("H5").Value = 36.93751
$ws.Range("I5").Value = 0.6498350963072504
$ws.Range("J5").Value = 0.6498350963072506
$ws.Range("M5").Value = 12.052907
$ws.Range("N5").Value = 36.158721
$ws.Range("O5").Value = 0.2349308184832226
$ws.Range("P5").Value = 0.2349308184832226
$ws.Range("Q5").Value = 148.4014576138567
$ws.Range("R5").Value = 1335.61311852471
$ws.Range("S5").Value = 0.1526662910545862
$ws.Range("T5").Value = 0.1526662910545862
$ws.Range("G6").Value = 12.31250333333333
$ws.Range("H6").Value = 36.93751
$ws.Range("I6").Value = 0.6498350963072504
$ws.Range("J6").Value = 0.6498350963072506
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1051546666666667
$ws.Range("N6").Value = 0.315464
$ws.Range("O6").Value = 0.002049635984690702
$ws.Range("P6").Value = 0.002049635984690701
$ws.Range("Q6").Value = 1.294717183848889
$ws.Range("R6").Value = 11.65245465464
$ws.Range("S6").Value = 0.001331925397506288
$ws.Range("T6").Value = 0.001331925397506288
$ws.Range("I7").Value = 0.3333514949915254
$ws.Range("J7").Value = 0.3333514949915254
$ws.Range("M7").Value = 4.889237666666667
$ws.Range("N7").Value = 14.667713
$ws.Range("O7").Value = 0.09529921759032918
$ws.Range("P7").Value = 0.09529921759032917
$ws.Range("Q7").Value = 30.88067282404178
$ws.Range("R7").Value = 277.926055416376
$ws.Range("S7").Value = 0.03176813665525891
$ws.Range("T7").Value = 0.03176813665525891
$ws.Range("I8").Value = 0.3333514949915254
$ws.Range("J8").Value = 0.3333514949915254
$ws.Range("O8").Value = 0.6664673019309815
$ws.Range("P8").Value = 0.6664673019309812
$ws.Range("S8").Value = 0.222167871461661
$ws.Range("T8").Value = 0.222167871461661
$ws.Range("I9").Value = 0.3333514949915254
$ws.Range("J9").Value = 0.3333514949915254
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.06428533333333333
$ws.Range("N9").Value = 0.192856
$ws.Range("O9").Value = 0.001253026010776221
$ws.Range("P9").Value = 0.001253026010776221
$ws.Range("Q9").Value = 0.4060294224568889
$ws.Range("R9").Value = 3.654264802112
$ws.Range("S9").Value = 0.0004176980939555206
$ws.Range("T9").Value = 0.0004176980939555206
$ws.Range("I10").Value = 0.3333514949915254
$ws.Range("J10").Value = 0.3333514949915254
$ws.Range("M10").Value = 12.052907
$ws.Range("N10").Value = 36.158721
$ws.Range("O10").Value = 0.2349308184832226
$ws.Range("P10").Value = 0.2349308184832226
$ws.Range("Q10").Value = 76.12677129262133
$ws.Range("R10").Value = 685.140941633592
$ws.Range("S10").Value = 0.07831453956096494
$ws.Range("T10").Value = 0.07831453956096496
$ws.Range("I11").Value = 0.3333514949915254
$ws.Range("J11").Value = 0.3333514949915254
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1051546666666667
$ws.Range("N11").Value = 0.315464
$ws.Range("O11").Value = 0.002049635984690702
$ws.Range("P11").Value = 0.002049635984690701
$ws.Range("Q11").Value = 0.6641622025031111
$ws.Range("R11").Value = 5.977459822528001
$ws.Range("S11").Value = 0.0006832492196850726
$ws.Range("T11").Value = 0.0006832492196850726
$ws.Range("G12").Value = 0.3185656666666667
$ws.Range("H12").Value = 0.955697
$ws.Range("I12").Value = 0.01681340870122405
$ws.Range("J12").Value = 0.01681340870122405
$ws.Range("M12").Value = 4.889237666666667
$ws.Range("N12").Value = 14.667713
$ws.Range("O12").Value = 0.09529921759032918
$ws.Range("P12").Value = 0.09529921759032917
$ws.Range("Q12").Value = 1.557543256773445
$ws.Range("R12").Value = 14.017889310961
$ws.Range("S12").Value = 0.001602304694253085
$ws.Range("T12").Value = 0.001602304694253084
$ws.Range("G13").Value = 0.3185656666666667
$ws.Range("H13").Value = 0.955697
$ws.Range("I13").Value = 0.01681340870122405
$ws.Range("J13").Value = 0.01681340870122405
$ws.Range("O13").Value = 0.6664673019309815
$ws.Range("P13").Value = 0.6664673019309812
$ws.Range("Q13").Value = 10.89255167282644
$ws.Range("R13").Value = 98.032965055438
$ws.Range("S13").Value = 0.01120558713336768
$ws.Range("T13").Value = 0.01120558713336768
$ws.Range("G14").Value = 0.3185656666666667
$ws.Range("H14").Value = 0.955697
$ws.Range("I14").Value = 0.01681340870122405
$ws.Range("J14").Value = 0.01681340870122405
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.06428533333333333
$ws.Range("N14").Value = 0.192856
$ws.Range("O14").Value = 0.001253026010776221
$ws.Range("P14").Value = 0.001253026010776221
$ws.Range("Q14").Value = 0.02047910007022222
$ws.Range("R14").Value = 0.184311900632
$ws.Range("S14").Value = 0.00002106763843244498
$ws.Range("T14").Value = 0.00002106763843244498
$ws.Range("G15").Value = 0.3185656666666667
$ws.Range("H15").Value = 0.955697
$ws.Range("I15").Value = 0.01681340870122405
$ws.Range("J15").Value = 0.01681340870122405
$ws.Range("M15").Value = 12.052907
$ws.Range("N15").Value = 36.158721
$ws.Range("O15").Value = 0.2349308184832226
$ws.Range("P15").Value = 0.2349308184832226
$ws.Range("Q15").Value = 3.839642353726334
$ws.Range("R15").Value = 34.556781183537
$ws.Range("S15").Value = 0.003949987867671503
$ws.Range("T15").Value = 0.003949987867671503
$ws.Range("G16").Value = 0.3185656666666667
$ws.Range("H16").Value = 0.955697
$ws.Range("I16").Value = 0.01681340870122405
$ws.Range("J16").Value = 0.01681340870122405
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1051546666666667
$ws.Range("N16").Value = 0.315464
$ws.Range("O16").Value = 0.002049635984690702
$ws.Range("P16").Value = 0.002049635984690701
$ws.Range("Q16").Value = 0.03349866648977778
$ws.Range("R16").Value = 0.301487998408
$ws.Range("S16").Value = 0.00003446136749934056
$ws.Range("T16").Value = 0.00003446136749934056
